{"js": "// The date line \"17/10/2019:\" becomes a date range \"17-24/10/2019:\"\n// (i.e. \"-24\" is inserted right after \"17\" and before \"/10/2019:\").\nconst body = context.document.body;\n\nconst fullHits = body.search(\"17/10/2019:\", { matchCase: true });\nfullHits.load(\"items\");\nawait context.sync();\n\nif (fullHits.items.length === 0) {\n  throw new Error(\"Could not find target text '17/10/2019:' in the document.\");\n}\n\n// There is exactly one paragraph containing \"17/10/2019:\" (the other date\n// occurrence in the document has no trailing colon), but scope the search to\n// the matched range regardless so we never touch an unrelated \"17\".\nconst target = fullHits.items[0];\n\nconst seventeenHits = target.search(\"17\", { matchCase: true });\nseventeenHits.load(\"items\");\nawait context.sync();\n\nconst seventeen = seventeenHits.items[0];\n\n// Insert \"-24\" immediately after \"17\", before \"/10/2019:\". insertText at the\n// collapsed point right after \"17\" inherits that run's formatting.\nseventeen.getRange(\"After\").insertText(\"-24\", \"Before\");\nawait context.sync();\n", "ps1": "# The date line \"17/10/2019:\" becomes a date range \"17-24/10/2019:\"\n# (i.e. \"-24\" is inserted right after \"17\" and before \"/10/2019:\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"17/10/2019:\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"17-24/10/2019:\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# 2 = wdReplaceOne: only replace the single (unique) match.\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
